$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C numeric value updates
$ws.Range("C2").Value = 44
$ws.Range("C3").Value = 41
$ws.Range("C5").Value = 37
$ws.Range("C6").Value = 37
$ws.Range("C8").Value = 42
$ws.Range("C10").Value = 41
$ws.Range("C11").Value = 41
$ws.Range("C12").Value = 40
$ws.Range("C13").Value = 40
$ws.Range("C15").Value = 37
$ws.Range("C16").Value = 44
$ws.Range("C17").Value = 43
$ws.Range("C18").Value = 35

# Column B text updates
$ws.Range("B4").Value = "<ereto>"
$ws.Range("B5").Value = "<three>"
$ws.Range("B8").Value = "<nom>"
$ws.Range("B9").Value = "<go>"
$ws.Range("B10").Value = "<and>"
$ws.Range("B11").Value = "<ey>"
$ws.Range("B14").Value = "<november>"
$ws.Range("B16").Value = "<mram>"
$ws.Range("B18").Value = "<nere>"
